# Trade #21 closed at 2026-02-17 04:08:20 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.02
$wsSummary.Range("B4").Value = 0.02
$wsSummary.Range("B5").Value = 0.02
$wsSummary.Range("B6").Value = 21
$wsSummary.Range("B7").Value = 7
$wsSummary.Range("B9").Value = 33.33

# --- Strategy Status sheet ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.02
$wsStatus.Range("D4").Value = 21
$wsStatus.Range("E4").Value = 0.02
$wsStatus.Range("F4").Value = 0.02
$wsStatus.Range("G4").Value = 33.33

# --- New trade row data, appended to both "All Trades" and "MarketMaking" sheets ---
function Add-TradeRow($ws) {
    $ws.Cells.Item(22, 1).Value = 21

    # Date column holds a plain text value like "2026-02-17" (not a real
    # Excel date) in this workbook, matching the existing rows above it.
    # Force text formatting first so COM doesn't auto-coerce the
    # date-shaped string into a serial date, then drop back to the default
    # (unstyled) cell style once the text value is committed.
    $dateCell = $ws.Cells.Item(22, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.Style = "Normal"

    $ws.Cells.Item(22, 3).Value = "04:08:14"
    $ws.Cells.Item(22, 4).Value = "MarketMaking"
    $ws.Cells.Item(22, 5).Value = "DOWN"
    $ws.Cells.Item(22, 6).Value = 0.58
    $ws.Cells.Item(22, 7).Value = 0.6
    $ws.Cells.Item(22, 8).Value = "CLOSED"
    $ws.Cells.Item(22, 9).Value = 3.4483
    $ws.Cells.Item(22, 10).Value = 0.02
    $ws.Cells.Item(22, 11).Value = 100.02
    $ws.Cells.Item(22, 12).Value = 0
    $ws.Cells.Item(22, 13).Value = 0
    $ws.Cells.Item(22, 14).Value = 0.6
    $ws.Cells.Item(22, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(22, 16).Value = "early_exit"
    $ws.Cells.Item(22, 17).Value = 0.12
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking
